# Updates cryptos list (price/volume refresh) to match the Sep 21 2023
# GitHub Actions data pull. Column D ("Price") holds numeric-looking text
# (e.g. "19.68"); assigning such a string straight to .Value lets the COM
# layer auto-coerce it to a Double, which would flip the cell from text to
# a number. Forcing NumberFormat="@" (Text) before the assignment keeps it
# a string, and ClearFormats() afterwards drops the now-unneeded explicit
# style so the cell's formatting matches the original (un-styled) cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.690.13'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -2.07%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.590.35'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -2.54%  '
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.16'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -2.31%  '
$ws.Range('E6').Value = '  -2.02%  '
$ws.Range('E8').Value = '  -2.58%  '
$ws.Range('E9').Value = '  -1.41%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.68'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -3.14%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0835'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -1.66%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.813.14'
$ws.Range('D12').ClearFormats()
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.574.78'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -3.45%  '
$ws.Range('E14').Value = '  -2.19%  '
$ws.Range('E15').Value = '  -2.93%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.61'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.80%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.646.77'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -2.07%  '
$ws.Range('E18').Value = '  -1.32%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '207.77'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -4.67%  '
$ws.Range('E20').Value = '  -0.14%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.76'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -2.54%  '
$ws.Range('E22').Value = '  -2.99%  '
$ws.Range('E23').Value = '  -2.95%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.90'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -2.07%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '147.23'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.56%  '
$ws.Range('E26').Value = '  -0.06%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.34'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +1.04%  '
$ws.Range('E28').Value = '  -3.56%  '
$ws.Range('E29').Value = '  -2.34%  '
$ws.Range('E30').Value = '  -0.53%  '
$ws.Range('E31').Value = '  -2.06%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.24'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -4.09%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.663'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +22.31%  '
$ws.Range('B34').Value = 'Maker'
$ws.Range('C34').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.323.57'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.39%  '
$ws.Range('B35').Value = 'InternetComputer(DFINITY)'
$ws.Range('C35').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.91'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -2.87%  '
$ws.Range('E36').Value = '  -4.06%  '
$ws.Range('E37').Value = '  -1.79%  '
$ws.Range('E38').Value = '  -1.54%  '
$ws.Range('E39').Value = '  -2.48%  '
$ws.Range('E41').Value = '  +3.62%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.789'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -1.19%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.17'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -3.87%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '63.47'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.43%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.726.06'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -2.37%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '89.96'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.89%  '
$ws.Range('E47').Value = '  -0.61%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.837'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +2.89%  '
$ws.Range('E49').Value = '  -0.79%  '
$ws.Range('E50').Value = '  +0.43%  '
